$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated values for column F (dSF) for specific rows, per the data repull.
$changes = @{
    "F2"  = -11
    "F4"  = -4
    "F5"  = -2
    "F7"  = -3
    "F8"  = -4
    "F9"  = -3
    "F11" = -3
    "F13" = -1
    "F15" = 0
    "F18" = -3
    "F19" = 1
    "F20" = -1
    "F25" = -4
    "F30" = 0
    "F32" = 3
    "F34" = -5
    "F36" = -6
    "F38" = 4
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
